# Update view-count values (column F) for a handful of events across the
# "展览" (sheet 1), "演出" (sheet 2) and "全部类型" (sheet 4) worksheets.
# These mirror a data refresh of the generated gh-pages content.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1918
$ws1.Range("F6").Value = 13458
$ws1.Range("F13").Value = 5
$ws1.Range("F16").Value = 2107
$ws1.Range("F23").Value = 293
$ws1.Range("F25").Value = 445
$ws1.Range("F27").Value = 37

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 137
$ws2.Range("F6").Value = 65
$ws2.Range("F8").Value = 571

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1918
$ws4.Range("F8").Value = 13458
$ws4.Range("F15").Value = 5
$ws4.Range("F20").Value = 2107
$ws4.Range("F25").Value = 137
$ws4.Range("F27").Value = 65
$ws4.Range("F31").Value = 293
$ws4.Range("F33").Value = 445
$ws4.Range("F36").Value = 571
$ws4.Range("F39").Value = 37
